$d = $word.ActiveDocument

$d.Content.Find.Execute("341×4=", $true, $false, $false, $false, $false, $true, 1, $false, "178×4=", 2) | Out-Null
$d.Content.Find.Execute("864×6=", $true, $false, $false, $false, $false, $true, 1, $false, "117×3=", 2) | Out-Null
$d.Content.Find.Execute("227×2=", $true, $false, $false, $false, $false, $true, 1, $false, "309×3=", 2) | Out-Null
$d.Content.Find.Execute("935×4=", $true, $false, $false, $false, $false, $true, 1, $false, "698×8=", 2) | Out-Null
$d.Content.Find.Execute("106×2=", $true, $false, $false, $false, $false, $true, 1, $false, "831×2=", 2) | Out-Null
$d.Content.Find.Execute("257×6=", $true, $false, $false, $false, $false, $true, 1, $false, "713×3=", 2) | Out-Null
$d.Content.Find.Execute("376×9=", $true, $false, $false, $false, $false, $true, 1, $false, "373×9=", 2) | Out-Null
$d.Content.Find.Execute("221×8=", $true, $false, $false, $false, $false, $true, 1, $false, "711×4=", 2) | Out-Null
$d.Content.Find.Execute("451×2=", $true, $false, $false, $false, $false, $true, 1, $false, "332×7=", 2) | Out-Null
$d.Content.Find.Execute("823×6=", $true, $false, $false, $false, $false, $true, 1, $false, "764×2=", 2) | Out-Null
$d.Content.Find.Execute("918×5=", $true, $false, $false, $false, $false, $true, 1, $false, "937×3=", 2) | Out-Null
$d.Content.Find.Execute("453×3=", $true, $false, $false, $false, $false, $true, 1, $false, "970×5=", 2) | Out-Null
$d.Content.Find.Execute("843×6=", $true, $false, $false, $false, $false, $true, 1, $false, "190×6=", 2) | Out-Null
$d.Content.Find.Execute("371×9=", $true, $false, $false, $false, $false, $true, 1, $false, "648×9=", 2) | Out-Null
$d.Content.Find.Execute("231×5=", $true, $false, $false, $false, $false, $true, 1, $false, "874×3=", 2) | Out-Null
$d.Content.Find.Execute("534×4=", $true, $false, $false, $false, $false, $true, 1, $false, "914×2=", 2) | Out-Null
$d.Content.Find.Execute("685×9=", $true, $false, $false, $false, $false, $true, 1, $false, "332×7=", 2) | Out-Null
$d.Content.Find.Execute("396×2=", $true, $false, $false, $false, $false, $true, 1, $false, "263×3=", 2) | Out-Null
$d.Content.Find.Execute("172×6=", $true, $false, $false, $false, $false, $true, 1, $false, "170×4=", 2) | Out-Null
$d.Content.Find.Execute("644×8=", $true, $false, $false, $false, $false, $true, 1, $false, "206×8=", 2) | Out-Null
$d.Content.Find.Execute("151×5=", $true, $false, $false, $false, $false, $true, 1, $false, "193×5=", 2) | Out-Null
$d.Content.Find.Execute("732×4=", $true, $false, $false, $false, $false, $true, 1, $false, "167×8=", 2) | Out-Null
$d.Content.Find.Execute("573×5=", $true, $false, $false, $false, $false, $true, 1, $false, "338×3=", 2) | Out-Null
$d.Content.Find.Execute("823×8=", $true, $false, $false, $false, $false, $true, 1, $false, "136×8=", 2) | Out-Null
$d.Content.Find.Execute("994×9=", $true, $false, $false, $false, $false, $true, 1, $false, "651×4=", 2) | Out-Null

Write-Host "Replacements complete"
